# Apply edit: add new row "FindKthLargest" to the "Heap" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Heap")
$ws.Activate()

$ws.Range("A7").Value = "FindKthLargest"

$ws.Range("B10").Select()
